# "Updated data for Dec 18" -- append new daily rows (117-122) to the
# clemsonUnivDaily sheet, correct a handful of previously-reported daily
# figures (rows 102-116), let the existing rolling-average / cumulative
# formulas fill down to the new rows, and refresh the weekly roll-up
# sheet (clemsonUnivWeekly) with the new week-30 summary row.

$wb  = $excel.ActiveWorkbook
$daily  = $wb.Worksheets.Item("clemsonUnivDaily")
$weekly = $wb.Worksheets.Item("clemsonUnivWeekly")

# ---------------------------------------------------------------------
# 1) Corrections to previously entered daily rows (101-116)
# ---------------------------------------------------------------------
$daily.Range("C102").Value = 392

$daily.Range("C103").Value = 131

$daily.Range("B104").Value = 3
$daily.Range("C104").Value = 636

$daily.Range("C105").Value = 1211

$daily.Range("C112").Value = 608

$daily.Range("C113").Value = 333

$daily.Range("B114").Value = 13
$daily.Range("C114").Value = 1308

$daily.Range("B115").Value = 14
$daily.Range("C115").Value = 545
$daily.Range("D115").Value = 9

$daily.Range("B116").Value = 17
$daily.Range("C116").Value = 724
$daily.Range("D116").Value = 13

# ---------------------------------------------------------------------
# 2) New daily rows 117-122 (formulas follow the same pattern used by
#    the rest of the table -- each column is a rolling calc off B:D)
# ---------------------------------------------------------------------
function Set-DailyRow {
    param(
        [int]$Row,
        [double]$B,
        [double]$C,
        [double]$D
    )
    $daily.Range("B$Row").Value = $B
    $daily.Range("C$Row").Value = $C
    $daily.Range("D$Row").Value = $D
    $daily.Range("E$Row").Formula = "=ABS(B$Row-D$Row)"
    $daily.Range("F$Row").Formula = "=AVERAGE(B$(($Row-6)):B$Row)"
    $daily.Range("G$Row").Formula = "=B$Row/C$Row"
    $daily.Range("H$Row").Formula = "=AVERAGE(G$(($Row-6)):G$Row)"
    $daily.Range("I$Row").Formula = "=(B$Row/26406)*100000"
    $daily.Range("J$Row").Formula = "=SUM(B$(($Row-9)):B$Row)"
    $daily.Range("K$Row").Formula = "=AVERAGE(J$(($Row-6)):J$Row)"
}

# Row 117 already has A117 (=12/18/20, serial 44177) filled in.
Set-DailyRow -Row 117 -B 0  -C 1   -D 0

# Row 118 already has A118 (serial 44178) filled in.
Set-DailyRow -Row 118 -B 1  -C 507 -D 0

# Row 119 already has A119 (serial 44179), L119=1 and M119=5 filled in.
Set-DailyRow -Row 119 -B 26 -C 590 -D 16

# Rows 120-126 are brand new rows.
$daily.Range("A120").Value = 44180
Set-DailyRow -Row 120 -B 11 -C 550 -D 9

$daily.Range("A121").Value = 44181
Set-DailyRow -Row 121 -B 14 -C 717 -D 8

$daily.Range("A122").Value = 44182
Set-DailyRow -Row 122 -B 0  -C 179 -D 0
$daily.Range("L122").Value = 1
$daily.Range("M122").Value = 4

$daily.Range("A123").Value = 44183
$daily.Range("A124").Value = 44184
$daily.Range("A125").Value = 44185
$daily.Range("A126").Value = 44186

# ---------------------------------------------------------------------
# 3) clemsonUnivWeekly: add the week-30 (12/18) summary row
# ---------------------------------------------------------------------
$weekly.Range("A30").Value = 44184
$weekly.Range("B30").Formula = "=SUM(clemsonUnivDaily!B118:B124)"
$weekly.Range("C30").Formula = "=(B30-B29)/B30"
$weekly.Range("D30").Formula = "=SUM(clemsonUnivDaily!C118:C124)"
$weekly.Range("E30").Formula = "=(D30-D29)/D30"
$weekly.Range("F30").Formula = "=SUM(clemsonUnivDaily!D118:D124)"
$weekly.Range("G30").Formula = "=SUM(clemsonUnivDaily!E118:E124)"
$weekly.Range("H30").Formula = "=B30/D30"
$weekly.Range("I30").Formula = "=(B30/26406)*100000"

try {
    $weekly.Range("B30").Errors.Item(1).Ignore = $true
} catch {
}

# ---------------------------------------------------------------------
# 4) View-state touch-ups from the diff (selection only -- topLeftCell /
#    window geometry are not represented after a COM round-trip)
# ---------------------------------------------------------------------
$daily.Range("A111:A117").Select()

$weekly.Activate()
$weekly.Range("B29").Select()

Write-Host "clemsonDashboard: applied Dec 18 data update"
